$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

function Replace-InRange($textRange, $oldStr, $newStr) {
    $full = $textRange.Text
    $idx = $full.IndexOf($oldStr)
    if ($idx -lt 0) {
        throw "Could not find string: [$oldStr]"
    }
    $sub = $textRange.Characters($idx + 1, $oldStr.Length)
    $sub.Text = $newStr
}

# Replacements are performed from the last occurrence in the text to the
# first, so that character offsets found earlier stay valid after each edit.

# 5) "Properties and outgoing edges are stored as "fat" graph properties."
#    -> "Properties/outgoing edges stored as "fat" graph properties."
$old5 = 'Properties and outgoing edges are stored as "fat" graph properties.'
$new5 = "Properties/outgoing edges stored as " + [char]0x201C + "fat" + [char]0x201D + " graph properties."
Replace-InRange $tr $old5 $new5

# 4) "Time dimension as first citizen;" -> "Time attribute as first citizen;"
$old4 = 'Time dimension as first citizen;'
$new4 = 'Time attribute as first citizen;'
Replace-InRange $tr $old4 $new4

# 3) "values > 8 bytes (e.g.<nbsp>strings, geometries) are stored in a dynamic
#     storage (RocksDB);"
#    -> "values > 8 bytes (e.g.<nbsp>strings, geometries) stored in dynamic
#        storage (RocksDB);"
# The non-breaking space (U+00A0) between "e.g." and "strings" cannot be
# reliably read back through TextRange.Text in this runtime (it round-trips
# as U+FFFD), so locate the span using the ASCII-safe prefix/suffix around it
# instead of matching the NBSP character itself.
$prefix3 = 'values > 8 bytes (e.g.'
$suffixOld3 = 'strings, geometries) are stored in a dynamic storage (RocksDB);'
$full = $tr.Text
$idxPrefix3 = $full.IndexOf($prefix3)
$idxSuffix3 = $full.IndexOf($suffixOld3)
if ($idxPrefix3 -lt 0 -or $idxSuffix3 -lt 0) {
    throw "Could not find NBSP-spanning paragraph"
}
$totalOldLen3 = ($idxSuffix3 + $suffixOld3.Length) - $idxPrefix3
$nbsp = [char]0x00A0
$new3 = "values > 8 bytes (e.g.$nbsp" + "strings, geometries) stored in dynamic storage (RocksDB);"
$sub3 = $tr.Characters($idxPrefix3 + 1, $totalOldLen3)
$sub3.Text = $new3

# 2) "Properties and edges are represented as a linked chain of pointers;"
#    -> "Properties/edges represented as linked chain of pointers;"
$old2 = 'Properties and edges are represented as a linked chain of pointers;'
$new2 = 'Properties/edges represented as linked chain of pointers;'
Replace-InRange $tr $old2 $new2

# 1) "Based on index-free adjacency through fixed-size records stored in
#     nodes, edges, and property files."
#    -> "Index-free adjacency through fixed-size records stored in nodes,
#        edges, and property files."
$old1 = 'Based on index-free adjacency through fixed-size records stored in nodes, edges, and property files.'
$new1 = 'Index-free adjacency through fixed-size records stored in nodes, edges, and property files.'
Replace-InRange $tr $old1 $new1
